$excel.Calculation = -4135
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.EnableCalculation = $false
$ws.Range("F11").Formula = "=_xll.MultiByteStrLen(REPT(C11,D11))"
